# Apply updated cryptocurrency Price (column D) and Volume(1h) (column E) figures to sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric (e.g. "293.93", "0.0701") but must stay plain text,
# exactly like the source data. A leading single-quote character forces Excel to
# store the value as text instead of reinterpreting it as a number.

$ws.Range("D2").Value = "'39.699.05"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").Value = "'2.218.95"

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'293.93"
$ws.Range("E5").Value = "  -5.33%  "

$ws.Range("D6").Value = "'83.86"
$ws.Range("E6").Value = "  -1.93%  "

$ws.Range("E7").Value = "  -2.75%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.465"
$ws.Range("E9").Value = "  -3.67%  "

$ws.Range("E10").Value = "  -3.36%  "

$ws.Range("D11").Value = "'29.66"
$ws.Range("E11").Value = "  -1.25%  "

$ws.Range("D12").Value = "'47.96"
$ws.Range("E12").Value = "  -8.47%  "

$ws.Range("E13").Value = "  -2.21%  "

$ws.Range("D14").Value = "'2.563.01"
$ws.Range("E14").Value = "  -5.38%  "

$ws.Range("D15").Value = "'6.28"
$ws.Range("E15").Value = "  -2.20%  "

$ws.Range("D16").Value = "'14.10"
$ws.Range("E16").Value = "  -4.28%  "

$ws.Range("D17").Value = "'2.220.81"
$ws.Range("E17").Value = "  -7.07%  "

$ws.Range("D18").Value = "'0.718"
$ws.Range("E18").Value = "  -5.36%  "

$ws.Range("D19").Value = "'39.620.06"
$ws.Range("E19").Value = "  -1.20%  "

$ws.Range("D20").Value = "'0.0₃0880"
$ws.Range("E20").Value = "  -2.45%  "

$ws.Range("E21").Value = "  -6.00%  "

$ws.Range("D22").Value = "'64.95"
$ws.Range("E22").Value = "  -4.54%  "

$ws.Range("D23").Value = "'10.41"
$ws.Range("E23").Value = "  -2.72%  "

$ws.Range("D24").Value = "'232.09"
$ws.Range("E24").Value = "  -1.33%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("E26").Value = "  -5.75%  "

$ws.Range("E27").Value = "  -0.46%  "

$ws.Range("D28").Value = "'22.76"
$ws.Range("E28").Value = "  -3.49%  "

$ws.Range("E29").Value = "  +2.48%  "

$ws.Range("D30").Value = "'9.17"

$ws.Range("D31").Value = "'32.10"
$ws.Range("E31").Value = "  -7.13%  "

$ws.Range("D32").Value = "'149.45"
$ws.Range("E32").Value = "  -2.62%  "

$ws.Range("E33").Value = "  -0.12%  "

$ws.Range("D34").Value = "'4.81"
$ws.Range("E34").Value = "  -5.75%  "

$ws.Range("D35").Value = "'2.38"
$ws.Range("E35").Value = "  -2.94%  "

$ws.Range("D36").Value = "'0.0701"
$ws.Range("E36").Value = "  -2.56%  "

$ws.Range("D37").Value = "'15.99"
$ws.Range("E37").Value = "  +2.99%  "

$ws.Range("D39").Value = "'0.0967"
$ws.Range("E39").Value = "  -1.44%  "

$ws.Range("E40").Value = "  -5.62%  "

$ws.Range("E41").Value = "  -4.03%  "

$ws.Range("E42").Value = "  -5.11%  "

$ws.Range("D43").Value = "'1.937.34"
$ws.Range("E43").Value = "  -1.42%  "

$ws.Range("D46").Value = "'9.40"
$ws.Range("E46").Value = "  +0.61%  "

$ws.Range("D47").Value = "'16.17"
$ws.Range("E47").Value = "  -8.20%  "

$ws.Range("D48").Value = "'2.59"
$ws.Range("E48").Value = "  -4.17%  "

$ws.Range("D49").Value = "'2.431.89"
$ws.Range("E49").Value = "  -5.37%  "

$ws.Range("D50").Value = "'70.59"
$ws.Range("E50").Value = "  +0.14%  "

$ws.Range("D51").Value = "'88.75"
$ws.Range("E51").Value = "  -4.66%  "
